# Drug resistance workbook: add "Drug sub-category", "Generation" and "Wave"
# columns between the existing "Drug category" and "Company code" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert three new columns (C, D, E). The existing "Company code"
#        column (old C) shifts to F. Excel inherits each new cell's style
#        from its left neighbour, which already reproduces almost all of
#        the target formatting (incl. the alternate Tahoma style on row 6
#        and the extra blank styled cell that appears in column G).
$ws.Range("C1:E12").EntireColumn.Insert() | Out-Null

# --- 2. Header row. Write new header text in the same order the strings
#        first appear in the finished workbook so they land on matching
#        shared-string indexes.
$ws.Range("C1").Value = "Drug sub-category"

# --- 3. "Drug sub-category" values for the NS5B_inhibitor drugs.
$ws.Range("C12").Value = "nucleotide_analogue"

$ws.Range("D1").Value = "Generation"
$ws.Range("E1").Value = "Wave"

$ws.Range("C11").Value = "non_nucleoside_palm_1_inhibitor"

# --- 4. "Generation" / "Wave" numeric values per drug row.
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 2

$ws.Range("D3").Value = 2
# E3 stays blank

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 2

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 2

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 2

$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1

$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 2

# Row 11 (dasabuvir) and row 12 (sofosbuvir) have no Generation/Wave value.

# --- 5. Row 12's new C/D/E cells use the alternate (Tahoma, themed-color)
#        style like the rest of that row, instead of the style inherited
#        from the insert. Copy the format from one of the already-correct
#        alternate-style cells (G2, created by the column insert) onto them.
$ws.Range("G2").Copy() | Out-Null
$ws.Range("C12:E12").PasteSpecial(-4122) | Out-Null
$ws.Range("C12").Value = "nucleotide_analogue"

# --- 6. Match the new column widths to the existing "Drug category" column.
$ws.Range("C1:E12").ColumnWidth = $ws.Range("B1").ColumnWidth

# --- 7. Restore the active selection shown in the finished workbook.
$ws.Range("D15").Select() | Out-Null
